$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44448
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 300000001
$ws.Cells.Item($row, 7).Value = "Rabanito"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 7000
$ws.Cells.Item($row, 12).Value = 7000
$ws.Cells.Item($row, 13).Value = 7000
$ws.Cells.Item($row, 14).Value = "$/docena de paquetes"
$ws.Cells.Item($row, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($row, 16).Value = 583
$ws.Cells.Item($row, 17).Value = 12
$ws.Cells.Item($row, 18).Value = "Hortaliza"
